$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the "Late" column (column N), matching the
# formatting that Excel copies in from the column on the left ("In Advance").
$leftWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab (was "Transactions").
$ws.Activate()
$ws.Range("J14").Select() | Out-Null
